$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated StatQuery text (column C, rows 2-4 share this string) ---
$newQuery = 'MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
WHERE demo.breed IN [''Boxer'']
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`'

$ws.Range("C2").Value = $newQuery
$ws.Range("C3").Value = $newQuery
$ws.Range("C4").Value = $newQuery

# --- Row heights (shrank because the replacement query text is shorter) ---
$ws.Rows.Item(2).RowHeight = 201.6
$ws.Rows.Item(3).RowHeight = 230.4
$ws.Rows.Item(4).RowHeight = 244.8

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 10
$ws.Columns.Item(2).ColumnWidth = 78.66666666666667
$ws.Columns.Item(3).ColumnWidth = 58
$ws.Columns.Item(4).ColumnWidth = 41.333333333333336
$ws.Columns.Item(5).ColumnWidth = 40.166666666666664

# --- View / selection: top-left scrolled to row 1, active cell now B4 ---
$ws.Activate()
$excel.Goto($ws.Range("B1"))
$ws.Range("B4").Select()
